# OpenData_Slovakia_Covid_DeathsCumulative.xlsx - append cumulative death
# counts for 2021-04-23 through 2021-05-24 (rows 190-221), matching the
# authoritative data refresh committed upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(190, 1).Value = 44309
$ws.Cells.Item(190, 2).Value = 11458
$ws.Cells.Item(190, 3).Value = 2197
$ws.Cells.Item(190, 4).Value = 13655
$ws.Cells.Item(191, 1).Value = 44310
$ws.Cells.Item(191, 2).Value = 11495
$ws.Cells.Item(191, 3).Value = 2210
$ws.Cells.Item(191, 4).Value = 13705
$ws.Cells.Item(192, 1).Value = 44311
$ws.Cells.Item(192, 2).Value = 11531
$ws.Cells.Item(192, 3).Value = 2227
$ws.Cells.Item(192, 4).Value = 13758
$ws.Cells.Item(193, 1).Value = 44312
$ws.Cells.Item(193, 2).Value = 11572
$ws.Cells.Item(193, 3).Value = 2236
$ws.Cells.Item(193, 4).Value = 13808
$ws.Cells.Item(194, 1).Value = 44313
$ws.Cells.Item(194, 2).Value = 11611
$ws.Cells.Item(194, 3).Value = 2247
$ws.Cells.Item(194, 4).Value = 13858
$ws.Cells.Item(195, 1).Value = 44314
$ws.Cells.Item(195, 2).Value = 11647
$ws.Cells.Item(195, 3).Value = 2255
$ws.Cells.Item(195, 4).Value = 13902
$ws.Cells.Item(196, 1).Value = 44315
$ws.Cells.Item(196, 2).Value = 11684
$ws.Cells.Item(196, 3).Value = 2266
$ws.Cells.Item(196, 4).Value = 13950
$ws.Cells.Item(197, 1).Value = 44316
$ws.Cells.Item(197, 2).Value = 11732
$ws.Cells.Item(197, 3).Value = 2273
$ws.Cells.Item(197, 4).Value = 14005
$ws.Cells.Item(198, 1).Value = 44317
$ws.Cells.Item(198, 2).Value = 11766
$ws.Cells.Item(198, 3).Value = 2285
$ws.Cells.Item(198, 4).Value = 14051
$ws.Cells.Item(199, 1).Value = 44318
$ws.Cells.Item(199, 2).Value = 11807
$ws.Cells.Item(199, 3).Value = 2292
$ws.Cells.Item(199, 4).Value = 14099
$ws.Cells.Item(200, 1).Value = 44319
$ws.Cells.Item(200, 2).Value = 11855
$ws.Cells.Item(200, 3).Value = 2302
$ws.Cells.Item(200, 4).Value = 14157
$ws.Cells.Item(201, 1).Value = 44320
$ws.Cells.Item(201, 2).Value = 11886
$ws.Cells.Item(201, 3).Value = 2316
$ws.Cells.Item(201, 4).Value = 14202
$ws.Cells.Item(202, 1).Value = 44321
$ws.Cells.Item(202, 2).Value = 11920
$ws.Cells.Item(202, 3).Value = 2324
$ws.Cells.Item(202, 4).Value = 14244
$ws.Cells.Item(203, 1).Value = 44322
$ws.Cells.Item(203, 2).Value = 11920
$ws.Cells.Item(203, 3).Value = 2335
$ws.Cells.Item(203, 4).Value = 14255
$ws.Cells.Item(204, 1).Value = 44323
$ws.Cells.Item(204, 2).Value = 11990
$ws.Cells.Item(204, 3).Value = 2339
$ws.Cells.Item(204, 4).Value = 14329
$ws.Cells.Item(205, 1).Value = 44324
$ws.Cells.Item(205, 2).Value = 12019
$ws.Cells.Item(205, 3).Value = 2346
$ws.Cells.Item(205, 4).Value = 14365
$ws.Cells.Item(206, 1).Value = 44325
$ws.Cells.Item(206, 2).Value = 12051
$ws.Cells.Item(206, 3).Value = 2351
$ws.Cells.Item(206, 4).Value = 14402
$ws.Cells.Item(207, 1).Value = 44326
$ws.Cells.Item(207, 2).Value = 12077
$ws.Cells.Item(207, 3).Value = 2359
$ws.Cells.Item(207, 4).Value = 14436
$ws.Cells.Item(208, 1).Value = 44327
$ws.Cells.Item(208, 2).Value = 12096
$ws.Cells.Item(208, 3).Value = 2370
$ws.Cells.Item(208, 4).Value = 14466
$ws.Cells.Item(209, 1).Value = 44328
$ws.Cells.Item(209, 2).Value = 12135
$ws.Cells.Item(209, 3).Value = 2380
$ws.Cells.Item(209, 4).Value = 14515
$ws.Cells.Item(210, 1).Value = 44329
$ws.Cells.Item(210, 2).Value = 12168
$ws.Cells.Item(210, 3).Value = 2387
$ws.Cells.Item(210, 4).Value = 14555
$ws.Cells.Item(211, 1).Value = 44330
$ws.Cells.Item(211, 2).Value = 12203
$ws.Cells.Item(211, 3).Value = 2392
$ws.Cells.Item(211, 4).Value = 14595
$ws.Cells.Item(212, 1).Value = 44331
$ws.Cells.Item(212, 2).Value = 12224
$ws.Cells.Item(212, 3).Value = 2395
$ws.Cells.Item(212, 4).Value = 14619
$ws.Cells.Item(213, 1).Value = 44332
$ws.Cells.Item(213, 2).Value = 12238
$ws.Cells.Item(213, 3).Value = 2402
$ws.Cells.Item(213, 4).Value = 14640
$ws.Cells.Item(214, 1).Value = 44333
$ws.Cells.Item(214, 2).Value = 12248
$ws.Cells.Item(214, 3).Value = 2403
$ws.Cells.Item(214, 4).Value = 14651
$ws.Cells.Item(215, 1).Value = 44334
$ws.Cells.Item(215, 2).Value = 12262
$ws.Cells.Item(215, 3).Value = 2405
$ws.Cells.Item(215, 4).Value = 14667
$ws.Cells.Item(216, 1).Value = 44335
$ws.Cells.Item(216, 2).Value = 12272
$ws.Cells.Item(216, 3).Value = 2405
$ws.Cells.Item(216, 4).Value = 14677
$ws.Cells.Item(217, 1).Value = 44336
$ws.Cells.Item(217, 2).Value = 12280
$ws.Cells.Item(217, 3).Value = 2406
$ws.Cells.Item(217, 4).Value = 14686
$ws.Cells.Item(218, 1).Value = 44337
$ws.Cells.Item(218, 2).Value = 12286
$ws.Cells.Item(218, 3).Value = 2406
$ws.Cells.Item(218, 4).Value = 14692
$ws.Cells.Item(219, 1).Value = 44338
$ws.Cells.Item(219, 2).Value = 12292
$ws.Cells.Item(219, 3).Value = 2406
$ws.Cells.Item(219, 4).Value = 14698
$ws.Cells.Item(220, 1).Value = 44339
$ws.Cells.Item(220, 2).Value = 12290
$ws.Cells.Item(220, 3).Value = 2407
$ws.Cells.Item(220, 4).Value = 14697
$ws.Cells.Item(221, 1).Value = 44340
$ws.Cells.Item(221, 2).Value = 12301
$ws.Cells.Item(221, 3).Value = 2407
$ws.Cells.Item(221, 4).Value = 14708

# Restore the default view: scroll back to the top-left and select A1,
# mirroring the saved workbook state (no stale topLeftCell/selection left
# over from row 189/143).
$ws.Range("A1").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
